$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for rows of 87b7ca0e.md (row 4) and
# e3d93a98.md (row 5) both currently show 2016-09-03 04:18:07 -> 04:18:57
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-03 04:18:57"
$wsOverview.Range("G5").Value = "2016-09-03 04:18:57"

# --- zh-cn sheet ---
# Status "ht" -> "mt" for 87b7ca0e (row 4) and e3d93a98 (row 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# Correspond Handoff Datetime 2016-09-03 04:17:57 -> 2016-09-03 04:18:53
$wsZhCn.Range("H4").Value = "2016-09-03 04:18:53"
$wsZhCn.Range("H5").Value = "2016-09-03 04:18:53"

# Correspond Handback DateTime 2016-09-03 04:18:27 -> 2016-09-03 04:19:15
$wsZhCn.Range("K4").Value = "2016-09-03 04:19:15"
$wsZhCn.Range("K5").Value = "2016-09-03 04:19:15"

# --- de-de sheet ---
# Status "ht" -> "mt" for 87b7ca0e (row 4) and e3d93a98 (row 5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# Correspond Handoff Datetime shares the same underlying text as the
# Overview sheet's "Latest HO Xliff Generate Date" (2016-09-03 04:18:07),
# which also moves to 2016-09-03 04:18:57
$wsDeDe.Range("H4").Value = "2016-09-03 04:18:57"
$wsDeDe.Range("H5").Value = "2016-09-03 04:18:57"

# Correspond Handback DateTime 2016-09-03 04:18:34 -> 2016-09-03 04:19:22
$wsDeDe.Range("K4").Value = "2016-09-03 04:19:22"
$wsDeDe.Range("K5").Value = "2016-09-03 04:19:22"
